# "fix bug exeded requeste in google drive"
#
# Bumps the price list's date stamp by one day and reprices the
# autoadhesive tope/zocalo SKUs (rows 29-33, column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds the list date as a date-formatted serial number (style 34,
# numFmtId 14). Read it back with Value2 (plain numeric, unlike Value
# which surfaces a date-flavoured COM Variant) and advance it one day.
$ws.Range("A1").Value = $ws.Range("A1").Value2 + 1

# Rows 29-32: price drops from 400 to 167.
$ws.Range("D29:D32").Value = 167

# Row 33: price drops from 1182 to 508.
$ws.Range("D33").Value = 508
